$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F, shifting old F (District) to G
$ws.Range("F1").EntireColumn.Insert()

# Populate the new Address column (F)
$ws.Range("F2").Value = "Address"
$ws.Range("F3").Value = "G J C Sringeri"
$ws.Range("F5").Value = "Govt. High School BemalkhedaHumnabad"
$ws.Range("F6").Value = "G H S Kundagol"
$ws.Range("F7").Value = "R L S High School"
$ws.Range("F8").Value = "G H S (RMSA) ParwatiBadami"
$ws.Range("F9").Value = "B T V P High School HallichedHumanbad"
$ws.Range("F10").Value = "Janatha Rural High School Lakshmipura"
$ws.Range("F11").Value = "S J F S High SchoolSanshiKundgol"
$ws.Range("F12").Value = "S V H S RamanathpuraArkalgudu"
$ws.Range("F13").Value = "Govt. Boys H S Humnabad"
$ws.Range("F14").Value = "S M H S KudlurTarikere"
$ws.Range("F16").Value = "G U H S AmingadHungund"
$ws.Range("F17").Value = "G H S HirehonnihalliKalaghatagi"
$ws.Range("F18").Value = "G J C (H S) AvinahalliSagar"
$ws.Range("F19").Value = "G H S Kuruvangi"
$ws.Range("F21").Value = "S R V P G H S KamadolliKundgol"
$ws.Range("F22").Value = "G H S HangarahalliHolenarsipura"
$ws.Range("F23").Value = "C K S Girls High SchoolK R Puram"
$ws.Range("F24").Value = "Adarsha Vidyalaya (RMSA) LokapurLaxanattiMudhol"
$ws.Range("F25").Value = "Govt. High School (Boys) NirnaHumnabad"
$ws.Range("F26").Value = "G J C (H S) ThatanahalliHolenarasipura"
$ws.Range("F27").Value = "Shreeveerbhadreshwar High School HalagalMudhol"
$ws.Range("F28").Value = "St. Joseph’s High School"
$ws.Range("F29").Value = "Sri A R High School KiggaSringeri"
$ws.Range("F30").Value = "S J G P U C Newtown Bhadravathi"
$ws.Range("F31").Value = "Y V M High School MorabNavalagund"
$ws.Range("F32").Value = "G H S DubalagundiHumnabad"
$ws.Range("F33").Value = "G H S TarlaghattaKundgol"
$ws.Range("F34").Value = "N V Kannya High School"
$ws.Range("F35").Value = "G H S Bannikuppe RamanagaraKanakapura"
$ws.Range("F36").Value = "R S H S Mailoor"
$ws.Range("F37").Value = "G J C (High School) HallimysoreHolenarasipura"
$ws.Range("F38").Value = "Malenadu High School"
$ws.Range("F40").Value = "Somaiya Vinay MandirHigh School SameerwadiMudhol"
$ws.Range("F41").Value = "Dayanand Hindi Vidyalay"
$ws.Range("F42").Value = "G H S S Madapura Kadur"
$ws.Range("F43").Value = "G H S Sangameshwarapete"
$ws.Range("F44").Value = "Govt. Urdu High School TippunagaraChannapatana"
$ws.Range("F45").Value = "G H S KarkeshwaraN R Pura"
$ws.Range("F46").Value = "Govt. H S MinajagiMuddebihal"
$ws.Range("F47").Value = "Sri Siddaganga Rural High SchoolDoddagangavadi"
$ws.Range("F48").Value = "G J C PaduvalahippeHolenarasipur"
$ws.Range("F49").Value = "Govt. Junior College Girls"
$ws.Range("F50").Value = "G H S BaradwadKundagol"
$ws.Range("F51").Value = "G H S MavanurHolenarasipur"
$ws.Range("F52").Value = "Girls High School RabkaviJamkhandi"
$ws.Range("F53").Value = "Govt. High School AlaghattaBirur"
$ws.Range("F54").Value = "Sheth Hanchaji Navalaji High School HudageriKundagola"
$ws.Range("F55").Value = "Govt. High School Aurad(B)"
$ws.Range("F56").Value = "P D J High School"
$ws.Range("F57").Value = "S R H S PallagatteJagalur"
$ws.Range("F58").Value = "G H S LingapuraHonnali"
